$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the threshold values for the rows that remain (alpha, beta, ratio)
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 5.6
$ws.Range("B4").Value = 0.7

# Remove the "theta_threshold_range" row entirely (row 5); this shifts the
# "pie_threshold_range" row (previously row 6) up into row 5.
$ws.Rows(5).Delete()

# Update the (now-shifted) pie_threshold_range Min/Max values.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 15

# Match the saved selection state.
$ws.Range("C12").Select()

# Match the saved page setup (paper size A4, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
